$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test row 3 ("Ürün arama çalışıyor mu?" / search-product test) currently reuses
# the "Kullanıcı giriş yapmış olmalı" precondition text in C4. Replace it with the
# correct precondition for the new search test, adding a fresh shared string.
$ws.Range("C4").Value = 'Kullanıcı anasayfa''da olmalı'

# Move the saved selection/cursor to C4 (the cell that was just edited).
$ws.Range("C4").Select() | Out-Null
